$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5's measurement columns (B:AH) to 2 decimal places (custom accuracy),
# matching the values recorded in the target workbook.
$ws.Range("B5").Value = 4.8
$ws.Range("C5").Value = 3.27
$ws.Range("F5").Value = 8.03
$ws.Range("H5").Value = 17.84
$ws.Range("I5").Value = 5.82
$ws.Range("J5").Value = 2.47
$ws.Range("K5").Value = 3.47
$ws.Range("L5").Value = 4.17
$ws.Range("M5").Value = 4.44
$ws.Range("N5").Value = 1.21
$ws.Range("P5").Value = 5.27
$ws.Range("Q5").Value = 3.4
$ws.Range("S5").Value = 0.37
$ws.Range("W5").Value = 6.95
$ws.Range("X5").Value = 3.74
$ws.Range("Y5").Value = 0.82
$ws.Range("AA5").Value = 3.07
$ws.Range("AC5").Value = 3.33
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 16.51
$ws.Range("AG5").Value = 1.83

# Drop the last data row (row 6) entirely -- "1000개" re-export trimmed the
# trailing sample, shrinking the sheet from A1:AH6 to A1:AH5.
$ws.Rows.Item(6).Delete()
